$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Translate the data-entry header row from Spanish/internal codes to the
# new English field names used by the updated database-saving code.
$ws.Range("B1").Value = "home"
$ws.Range("C1").Value = "date"
$ws.Range("D1").Value = "reason"
$ws.Range("E1").Value = "place"
$ws.Range("H1").Value = "general_service"
$ws.Range("I1").Value = "general_donor"
$ws.Range("J1").Value = "beneficiaries_served"
$ws.Range("K1").Value = "individual_service"
$ws.Range("L1").Value = "individual_donor"

# The longer English headers in K (individual_service) and L
# (individual_donor) no longer fit the old shared default width, so widen
# those two columns to fit their contents (mirrors an Excel AutoFit of K:L).
$ws.Columns.Item(11).ColumnWidth = 18.666666666666664
$ws.Columns.Item(12).ColumnWidth = 17.666666666666664

# Restore the active selection the author left on the sheet before saving.
$ws.Range("B4").Select()
